# Scheduled runner update: refresh market-board derived pricing/profit
# figures (currentAveragePrice*, LevePrice*, LeveProfit*) across the
# ALC/ARM/BSM/CRP/CUL/LTW/WVR sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 40000428
$ws.Range("I33").Value = 457.52173
$ws.Range("J33").Value = 500000100
$ws.Range("K33").Value = 457.52173
$ws.Range("L33").Value = 500000100
$ws.Range("M33").Value = -228.52173
$ws.Range("N33").Value = -500000558

$ws.Range("H76").Value = 3088.889
$ws.Range("I76").Value = 3114.2856
$ws.Range("K76").Value = 3114.2856
$ws.Range("M76").Value = -2799.2856

$ws.Range("H79").Value = 3088.889
$ws.Range("I79").Value = 3114.2856
$ws.Range("K79").Value = 3114.2856
$ws.Range("M79").Value = -2022.2856

$ws.Range("H99").Value = 1267.6316
$ws.Range("I99").Value = 853.6
$ws.Range("J99").Value = 2820.25
$ws.Range("K99").Value = 2560.8
$ws.Range("L99").Value = 8460.75
$ws.Range("M99").Value = -1062.8
$ws.Range("N99").Value = -11456.75

$ws.Range("H113").Value = 5124.5
$ws.Range("I113").Value = 3115
$ws.Range("J113").Value = 5985.7144
$ws.Range("K113").Value = 3115
$ws.Range("L113").Value = 5985.7144
$ws.Range("M113").Value = 139
$ws.Range("N113").Value = -12493.7144

$ws.Range("H129").Value = 973.34424
$ws.Range("I129").Value = 663.4
$ws.Range("J129").Value = 1034.1177
$ws.Range("K129").Value = 1990.2
$ws.Range("L129").Value = 3102.3531
$ws.Range("M129").Value = 3009.8
$ws.Range("N129").Value = -13102.3531

$ws.Range("H137").Value = 2780364.2
$ws.Range("I137").Value = 2237.2856
$ws.Range("K137").Value = 6711.8568
$ws.Range("M137").Value = -4161.8568

$ws.Range("H139").Value = 44029.8
$ws.Range("J139").Value = 44029.8
$ws.Range("L139").Value = 44029.8
$ws.Range("N139").Value = -54309.8

$ws.Range("H140").Value = 40480
$ws.Range("J140").Value = 40480
$ws.Range("L140").Value = 40480
$ws.Range("N140").Value = -50840

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2986.9714
$ws.Range("I32").Value = 2689.0715
$ws.Range("J32").Value = 4178.5713
$ws.Range("K32").Value = 2689.0715
$ws.Range("L32").Value = 4178.5713
$ws.Range("M32").Value = -2402.0715
$ws.Range("N32").Value = -4752.5713

$ws.Range("H122").Value = 1943.25
$ws.Range("I122").Value = 1690.1818
$ws.Range("K122").Value = 5070.5454
$ws.Range("M122").Value = -2620.5454

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H59").Value = 46792.4
$ws.Range("J59").Value = 46792.4
$ws.Range("L59").Value = 46792.4
$ws.Range("N59").Value = -48486.4

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H60").Value = 11056.917
$ws.Range("I60").Value = 10000
$ws.Range("J60").Value = 11102.869
$ws.Range("K60").Value = 10000
$ws.Range("L60").Value = 11102.869
$ws.Range("M60").Value = -9489
$ws.Range("N60").Value = -12124.869

$ws.Range("H110").Value = 30000
$ws.Range("J110").Value = 30000
$ws.Range("L110").Value = 30000
$ws.Range("N110").Value = -38180

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 588.06665
$ws.Range("I5").Value = 416.09525
$ws.Range("J5").Value = 989.3333
$ws.Range("K5").Value = 1248.28575
$ws.Range("L5").Value = 2967.9999
$ws.Range("M5").Value = -1136.28575
$ws.Range("N5").Value = -3191.9999

$ws.Range("H17").Value = 1543.3684
$ws.Range("I17").Value = 255.53847
$ws.Range("J17").Value = 4333.6665
$ws.Range("K17").Value = 766.61541
$ws.Range("L17").Value = 13000.9995
$ws.Range("M17").Value = -597.61541
$ws.Range("N17").Value = -13338.9995

$ws.Range("H95").Value = 9284.333000000001
$ws.Range("J95").Value = 9284.333000000001
$ws.Range("L95").Value = 27852.999
$ws.Range("N95").Value = -31970.999

$ws.Range("H113").Value = 643.9048
$ws.Range("I113").Value = 653.9
$ws.Range("J113").Value = 634.8182
$ws.Range("K113").Value = 1961.7
$ws.Range("L113").Value = 1904.4546
$ws.Range("M113").Value = 208.3000000000002
$ws.Range("N113").Value = -6244.4546

$ws.Range("H122").Value = 1013.087
$ws.Range("I122").Value = 536.0909
$ws.Range("K122").Value = 4824.8181
$ws.Range("M122").Value = -2374.8181

$ws.Range("H131").Value = 873.8788
$ws.Range("I131").Value = 470.57144
$ws.Range("J131").Value = 904.56525
$ws.Range("K131").Value = 1411.71432
$ws.Range("L131").Value = 2713.69575
$ws.Range("M131").Value = 3628.28568
$ws.Range("N131").Value = -12793.69575

$ws.Range("H135").Value = 588.06665
$ws.Range("I135").Value = 416.09525
$ws.Range("J135").Value = 989.3333
$ws.Range("K135").Value = 3744.85725
$ws.Range("L135").Value = 8903.9997
$ws.Range("M135").Value = -1209.85725
$ws.Range("N135").Value = -13973.9997

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3443.2144
$ws.Range("I7").Value = 3285.7144
$ws.Range("J7").Value = 3600.7144
$ws.Range("K7").Value = 3285.7144
$ws.Range("L7").Value = 3600.7144
$ws.Range("M7").Value = -3173.7144
$ws.Range("N7").Value = -3824.7144

$ws.Range("H40").Value = 2389.8
$ws.Range("I40").Value = 2389.8
$ws.Range("K40").Value = 2389.8
$ws.Range("M40").Value = -2253.8

$ws.Range("H122").Value = 4542.5713
$ws.Range("I122").Value = 4000
$ws.Range("J122").Value = 4949.5
$ws.Range("K122").Value = 12000
$ws.Range("L122").Value = 14848.5
$ws.Range("M122").Value = -9550
$ws.Range("N122").Value = -19748.5

$ws.Range("H126").Value = 3443.2144
$ws.Range("I126").Value = 3285.7144
$ws.Range("J126").Value = 3600.7144
$ws.Range("K126").Value = 9857.143199999999
$ws.Range("L126").Value = 10802.1432
$ws.Range("M126").Value = -7387.143199999999
$ws.Range("N126").Value = -15742.1432

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 153370.62
$ws.Range("I122").Value = 4240
$ws.Range("J122").Value = 302501.25
$ws.Range("K122").Value = 12720
$ws.Range("L122").Value = 907503.75
$ws.Range("M122").Value = -10270
$ws.Range("N122").Value = -912403.75

$ws.Range("H126").Value = 1630.85
$ws.Range("I126").Value = 1080.8
$ws.Range("K126").Value = 3242.4
$ws.Range("M126").Value = -772.3999999999996

$ws.Range("H136").Value = 3702.16
$ws.Range("I136").Value = 4262.0586
$ws.Range("J136").Value = 2512.375
$ws.Range("K136").Value = 12786.1758
$ws.Range("L136").Value = 7537.125
$ws.Range("M136").Value = -10236.1758
$ws.Range("N136").Value = -12637.125
